$d = $word.ActiveDocument

# Build the WordprocessingML fragment for the three new paragraphs that
# get appended after "The Contributor will write in "Vivaldi" font as a
# reply.":
#   1. An empty paragraph (sz/szCs 46, en-IN)
#   2. "Hi, this is Bhutu Banik and I have made the changes as per your
#      request." in Vivaldi font (sz/szCs 46, en-IN), with proofErr spell
#      markers bracketing "Bhutu"
#   3. An empty paragraph (sz/szCs 30, en-IN) matching the document's
#      normal body paragraph mark formatting

$body = '<w:p><w:pPr><w:rPr><w:sz w:val="46"/><w:szCs w:val="46"/><w:lang w:val="en-IN"/></w:rPr></w:pPr></w:p>' + `
  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Vivaldi" w:hAnsi="Vivaldi"/><w:sz w:val="46"/><w:szCs w:val="46"/><w:lang w:val="en-IN"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Vivaldi" w:hAnsi="Vivaldi"/><w:sz w:val="46"/><w:szCs w:val="46"/><w:lang w:val="en-IN"/></w:rPr><w:t xml:space="preserve">Hi, this is </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Vivaldi" w:hAnsi="Vivaldi"/><w:sz w:val="46"/><w:szCs w:val="46"/><w:lang w:val="en-IN"/></w:rPr><w:t>Bhutu</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Vivaldi" w:hAnsi="Vivaldi"/><w:sz w:val="46"/><w:szCs w:val="46"/><w:lang w:val="en-IN"/></w:rPr><w:t xml:space="preserve"> Banik and I have made the changes as per your request.</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="en-IN"/></w:rPr></w:pPr></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document></pkg:xmlData>' + `
  '</pkg:part></pkg:package>'

# Insert at a fresh zero-length Range anchored to the very end of the
# document's story so the new paragraphs are appended after the last
# existing paragraph, without disturbing that paragraph's own run/text.
$endPos = $d.Content.End
$r = $d.Range($endPos, $endPos)
$r.InsertXML($xml)
